# Update NATMI ligand-receptor pair values (Chad-Itga2) with recomputed TPM-based statistics.
# For each "Sending cluster" group the underlying ligand/receptor expression values were
# recomputed from the refreshed TPM matrix, which cascades into the detection rate,
# average/total expression, derived-specificity and edge-weight columns (F through T).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04253433333333333
$ws.Range("H2").Value = 0.127603
$ws.Range("I2").Value = 0.01491315744324368
$ws.Range("J2").Value = 0.01491315744324368
$ws.Range("M2").Value = 6.066605666666667
$ws.Range("N2").Value = 18.199817
$ws.Range("O2").Value = 0.8497846287916651
$ws.Range("P2").Value = 0.8497846287916652
$ws.Range("Q2").Value = 0.2580390276278889
$ws.Range("R2").Value = 2.322351248651
$ws.Range("S2").Value = 0.01267297196201849
$ws.Range("T2").Value = 0.01267297196201849

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04253433333333333
$ws.Range("H3").Value = 0.127603
$ws.Range("I3").Value = 0.01491315744324368
$ws.Range("J3").Value = 0.01491315744324368
$ws.Range("O3").Value = 0.1196497582104962
$ws.Range("P3").Value = 0.1196497582104962
$ws.Range("Q3").Value = 0.03633192013422221
$ws.Range("R3").Value = 0.326987281208
$ws.Range("S3").Value = 0.001784355682239168
$ws.Range("T3").Value = 0.001784355682239168

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04253433333333333
$ws.Range("H4").Value = 0.127603
$ws.Range("I4").Value = 0.01491315744324368
$ws.Range("J4").Value = 0.01491315744324368
$ws.Range("M4").Value = 0.1824346666666667
$ws.Range("N4").Value = 0.547304
$ws.Range("O4").Value = 0.02555468148257719
$ws.Range("P4").Value = 0.02555468148257719
$ws.Range("Q4").Value = 0.007759736923555556
$ws.Range("R4").Value = 0.06983763231199999
$ws.Range("S4").Value = 0.0003811009883616174
$ws.Range("T4").Value = 0.0003811009883616174

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.04253433333333333
$ws.Range("H5").Value = 0.127603
$ws.Range("I5").Value = 0.01491315744324368
$ws.Range("J5").Value = 0.01491315744324368
$ws.Range("M5").Value = 0.035773
$ws.Range("N5").Value = 0.107319
$ws.Range("O5").Value = 0.005010931515261538
$ws.Range("P5").Value = 0.005010931515261539
$ws.Range("Q5").Value = 0.001521580706333333
$ws.Range("R5").Value = 0.013694226357
$ws.Range("S5").Value = 0.00007472881062440695
$ws.Range("T5").Value = 0.00007472881062440695

# Row 6
$ws.Range("I6").Value = 0.2847488267267417
$ws.Range("J6").Value = 0.2847488267267417
$ws.Range("M6").Value = 6.066605666666667
$ws.Range("N6").Value = 18.199817
$ws.Range("O6").Value = 0.8497846287916651
$ws.Range("P6").Value = 0.8497846287916652
$ws.Range("Q6").Value = 4.926945259338
$ws.Range("R6").Value = 44.342507334042
$ws.Range("S6").Value = 0.2419751760188464
$ws.Range("T6").Value = 0.2419751760188464

# Row 7
$ws.Range("I7").Value = 0.2847488267267417
$ws.Range("J7").Value = 0.2847488267267417
$ws.Range("O7").Value = 0.1196497582104962
$ws.Range("P7").Value = 0.1196497582104962
$ws.Range("S7").Value = 0.03407012826857712
$ws.Range("T7").Value = 0.03407012826857712

# Row 8
$ws.Range("I8").Value = 0.2847488267267417
$ws.Range("J8").Value = 0.2847488267267417
$ws.Range("M8").Value = 0.1824346666666667
$ws.Range("N8").Value = 0.547304
$ws.Range("O8").Value = 0.02555468148257719
$ws.Range("P8").Value = 0.02555468148257719
$ws.Range("Q8").Value = 0.148162855056
$ws.Range("R8").Value = 1.333465695504
$ws.Range("S8").Value = 0.007276665569539447
$ws.Range("T8").Value = 0.007276665569539448

# Row 9
$ws.Range("I9").Value = 0.2847488267267417
$ws.Range("J9").Value = 0.2847488267267417
$ws.Range("M9").Value = 0.035773
$ws.Range("N9").Value = 0.107319
$ws.Range("O9").Value = 0.005010931515261538
$ws.Range("P9").Value = 0.005010931515261539
$ws.Range("Q9").Value = 0.029052755766
$ws.Range("R9").Value = 0.261474801894
$ws.Range("S9").Value = 0.001426856869778777
$ws.Range("T9").Value = 0.001426856869778777

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 1.236598666666667
$ws.Range("H10").Value = 3.709796
$ws.Range("I10").Value = 0.4335695228977033
$ws.Range("J10").Value = 0.4335695228977032
$ws.Range("M10").Value = 6.066605666666667
$ws.Range("N10").Value = 18.199817
$ws.Range("O10").Value = 0.8497846287916651
$ws.Range("P10").Value = 0.8497846287916652
$ws.Range("Q10").Value = 7.501956478592445
$ws.Range("R10").Value = 67.517608307332
$ws.Range("S10").Value = 0.3684407160710041
$ws.Range("T10").Value = 0.3684407160710041

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 1.236598666666667
$ws.Range("H11").Value = 3.709796
$ws.Range("I11").Value = 0.4335695228977033
$ws.Range("J11").Value = 0.4335695228977032
$ws.Range("O11").Value = 0.1196497582104962
$ws.Range("P11").Value = 0.1196497582104962
$ws.Range("Q11").Value = 1.056276200295111
$ws.Range("R11").Value = 9.506485802655998
$ws.Range("S11").Value = 0.05187648858215038
$ws.Range("T11").Value = 0.05187648858215039

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 1.236598666666667
$ws.Range("H12").Value = 3.709796
$ws.Range("I12").Value = 0.4335695228977033
$ws.Range("J12").Value = 0.4335695228977032
$ws.Range("M12").Value = 0.1824346666666667
$ws.Range("N12").Value = 0.547304
$ws.Range("O12").Value = 0.02555468148257719
$ws.Range("P12").Value = 0.02555468148257719
$ws.Range("Q12").Value = 0.2255984655537778
$ws.Range("R12").Value = 2.030386189984
$ws.Range("S12").Value = 0.01107973105820376
$ws.Range("T12").Value = 0.01107973105820376

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 1.236598666666667
$ws.Range("H13").Value = 3.709796
$ws.Range("I13").Value = 0.4335695228977033
$ws.Range("J13").Value = 0.4335695228977032
$ws.Range("M13").Value = 0.035773
$ws.Range("N13").Value = 0.107319
$ws.Range("O13").Value = 0.005010931515261538
$ws.Range("P13").Value = 0.005010931515261539
$ws.Range("Q13").Value = 0.04423684410266666
$ws.Range("R13").Value = 0.398131596924
$ws.Range("S13").Value = 0.00217258718634501
$ws.Range("T13").Value = 0.002172587186345011

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.03742733333333333
$ws.Range("H14").Value = 0.112282
$ws.Range("I14").Value = 0.01312256878006228
$ws.Range("J14").Value = 0.01312256878006228
$ws.Range("M14").Value = 6.066605666666667
$ws.Range("N14").Value = 18.199817
$ws.Range("O14").Value = 0.8497846287916651
$ws.Range("P14").Value = 0.8497846287916652
$ws.Range("Q14").Value = 0.2270568724882222
$ws.Range("R14").Value = 2.043511852394
$ws.Range("S14").Value = 0.01115135723955832
$ws.Range("T14").Value = 0.01115135723955832

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.03742733333333333
$ws.Range("H15").Value = 0.112282
$ws.Range("I15").Value = 0.01312256878006228
$ws.Range("J15").Value = 0.01312256878006228
$ws.Range("O15").Value = 0.1196497582104962
$ws.Range("P15").Value = 0.1196497582104962
$ws.Range("Q15").Value = 0.03196962968355555
$ws.Range("R15").Value = 0.287726667152
$ws.Range("S15").Value = 0.001570112181635057
$ws.Range("T15").Value = 0.001570112181635058

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.03742733333333333
$ws.Range("H16").Value = 0.112282
$ws.Range("I16").Value = 0.01312256878006228
$ws.Range("J16").Value = 0.01312256878006228
$ws.Range("M16").Value = 0.1824346666666667
$ws.Range("N16").Value = 0.547304
$ws.Range("O16").Value = 0.02555468148257719
$ws.Range("P16").Value = 0.02555468148257719
$ws.Range("Q16").Value = 0.006828043080888889
$ws.Range("R16").Value = 0.061452387728
$ws.Range("S16").Value = 0.000335343065407703
$ws.Range("T16").Value = 0.000335343065407703

# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.03742733333333333
$ws.Range("H17").Value = 0.112282
$ws.Range("I17").Value = 0.01312256878006228
$ws.Range("J17").Value = 0.01312256878006228
$ws.Range("M17").Value = 0.035773
$ws.Range("N17").Value = 0.107319
$ws.Range("O17").Value = 0.005010931515261538
$ws.Range("P17").Value = 0.005010931515261539
$ws.Range("Q17").Value = 0.001338887995333333
$ws.Range("R17").Value = 0.012049991958
$ws.Range("S17").Value = 0.00006575629346120122
$ws.Range("T17").Value = 0.00006575629346120123

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.7234323333333333
$ws.Range("H18").Value = 2.170297
$ws.Range("I18").Value = 0.253645924152249
$ws.Range("J18").Value = 0.253645924152249
$ws.Range("M18").Value = 6.066605666666667
$ws.Range("N18").Value = 18.199817
$ws.Range("O18").Value = 0.8497846287916651
$ws.Range("P18").Value = 0.8497846287916652
$ws.Range("Q18").Value = 4.388778692849889
$ws.Range("R18").Value = 39.499008235649
$ws.Range("S18").Value = 0.2155444075002378
$ws.Range("T18").Value = 0.2155444075002378

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.7234323333333333
$ws.Range("H19").Value = 2.170297
$ws.Range("I19").Value = 0.253645924152249
$ws.Range("J19").Value = 0.253645924152249
$ws.Range("O19").Value = 0.1196497582104962
$ws.Range("P19").Value = 0.1196497582104962
$ws.Range("Q19").Value = 0.6179404659102221
$ws.Range("R19").Value = 5.561464193192
$ws.Range("S19").Value = 0.03034867349589445
$ws.Range("T19").Value = 0.03034867349589445

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.7234323333333333
$ws.Range("H20").Value = 2.170297
$ws.Range("I20").Value = 0.253645924152249
$ws.Range("J20").Value = 0.253645924152249
$ws.Range("M20").Value = 0.1824346666666667
$ws.Range("N20").Value = 0.547304
$ws.Range("O20").Value = 0.02555468148257719
$ws.Range("P20").Value = 0.02555468148257719
$ws.Range("Q20").Value = 0.1319791365875556
$ws.Range("R20").Value = 1.187812229288
$ws.Range("S20").Value = 0.006481840801064656
$ws.Range("T20").Value = 0.006481840801064656

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.7234323333333333
$ws.Range("H21").Value = 2.170297
$ws.Range("I21").Value = 0.253645924152249
$ws.Range("J21").Value = 0.253645924152249
$ws.Range("M21").Value = 0.035773
$ws.Range("N21").Value = 0.107319
$ws.Range("O21").Value = 0.005010931515261538
$ws.Range("P21").Value = 0.005010931515261539
$ws.Range("Q21").Value = 0.02587934486033333
$ws.Range("R21").Value = 0.232914103743
$ws.Range("S21").Value = 0.001271002355052142
$ws.Range("T21").Value = 0.001271002355052142

